# Update statistics cells produced by a re-run of the capri bootstrap
# computation (BIC sheet "LUAD-bic" and AIC sheet "LUAD-aic").
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LUAD-bic")
$ws2 = $wb.Worksheets.Item("LUAD-aic")

# --- LUAD-bic sheet -------------------------------------------------
$ws1.Range("K2").Value = 0.14375
$ws1.Range("L2").Value = 0.04218428485690955
$ws1.Range("N2").Value = 0.0416666666666667
$ws1.Range("D3").Value = 2.0
$ws1.Range("F3").Value = [double]"5.11067883527512E-4"
$ws1.Range("G3").Value = [double]"4.21328460987917E-8"
$ws1.Range("K3").Value = 0.125
$ws1.Range("L3").Value = 0.0
$ws1.Range("M3").Value = 0.125
$ws1.Range("N3").Value = 0.0
$ws1.Range("I4").Value = 40.0
$ws1.Range("K4").Value = 0.125
$ws1.Range("L4").Value = 0.05103103630798288
$ws1.Range("M4").Value = 0.13125
$ws1.Range("N4").Value = 0.0354778882623467
$ws1.Range("D5").Value = 3.0
$ws1.Range("F5").Value = [double]"4.23751999471512E-10"
$ws1.Range("G5").Value = [double]"2.46328156654618E-10"
$ws1.Range("K5").Value = 0.1375
$ws1.Range("L5").Value = 0.03952847075210474
$ws1.Range("M5").Value = 0.10625
$ws1.Range("N5").Value = 0.030190368221228
$ws1.Range("I6").Value = 30.0
$ws1.Range("K6").Value = 0.125
$ws1.Range("L6").Value = 0.0
$ws1.Range("M6").Value = 0.11875
$ws1.Range("N6").Value = 0.0197642353760524
$ws1.Range("I7").Value = 20.0
$ws1.Range("K24").Value = 0.5375
$ws1.Range("L24").Value = 0.07905694150420949

# --- LUAD-aic sheet -------------------------------------------------
$ws2.Range("K2").Value = 0.14375
$ws2.Range("L2").Value = 0.04218428485690955
$ws2.Range("M2").Value = 0.15
$ws2.Range("N2").Value = 0.0437003686737563
$ws2.Range("D3").Value = 2.0
$ws2.Range("F3").Value = [double]"5.11067883527512E-4"
$ws2.Range("G3").Value = [double]"4.21328460987917E-8"
$ws2.Range("K3").Value = 0.125
$ws2.Range("L3").Value = 0.0
$ws2.Range("M3").Value = 0.125
$ws2.Range("N3").Value = 0.0
$ws2.Range("I4").Value = 40.0
$ws2.Range("K4").Value = 0.125
$ws2.Range("L4").Value = 0.05103103630798288
$ws2.Range("M4").Value = 0.125
$ws2.Range("N4").Value = 0.0416666666666667
$ws2.Range("D5").Value = 3.0
$ws2.Range("F5").Value = [double]"4.23751999471512E-10"
$ws2.Range("G5").Value = [double]"2.46328156654618E-10"
$ws2.Range("K5").Value = 0.1375
$ws2.Range("L5").Value = 0.03952847075210474
$ws2.Range("M5").Value = 0.10625
$ws2.Range("N5").Value = 0.030190368221228
$ws2.Range("I6").Value = 30.0
$ws2.Range("K6").Value = 0.125
$ws2.Range("L6").Value = 0.0
$ws2.Range("M6").Value = 0.11875
$ws2.Range("N6").Value = 0.0197642353760524
$ws2.Range("I7").Value = 20.0
$ws2.Range("M10").Value = 0.1875
$ws2.Range("N10").Value = 0.0
$ws2.Range("K25").Value = 0.5375
$ws2.Range("L25").Value = 0.07905694150420949
$ws2.Range("A3").NumberFormat = "@"
$ws2.Range("A3").Value = "8"
$ws2.Range("A3").ClearFormats()
$ws2.Range("A5").NumberFormat = "@"
$ws2.Range("A5").Value = "17"
$ws2.Range("A5").ClearFormats()
